# Update "想去人数" (want-to-go count) figures in column F
# for the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 297
    3  = 301
    5  = 28
    6  = 320
    7  = 9933
    11 = 128
    13 = 44
    14 = 31
    16 = 29
    17 = 285
    18 = 800
    19 = 48
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
